$wb = $excel.ActiveWorkbook

$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBS = $wb.Worksheets.Item("DBS")

# Add a new row of data to the DBS sheet (row 6)
$wsDBS.Range("A6").Value = "findCustNoFirst"
$wsDBS.Range("B6").Value = "CustNo = ,AND FacmNo >= ,AND FacmNo <= ,AND BormNo >= ,AND BormNo <= ,AND RepayDate >= ,AND RepayDate <="
$wsDBS.Range("C6").Value = "LogNo DESC"

# Update the selection on DBS sheet to A7
$wsDBS.Range("A7").Select()

# Make DBS the active (selected) sheet/tab
$wsDBS.Activate()
